$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.04244391055416
$ws.Cells.Item(2, 4).Value = 1.049465110399412
$ws.Cells.Item(2, 5).Value = 1.04955725944453
$ws.Cells.Item(2, 6).Value = 1.057911667164067
$ws.Cells.Item(2, 9).Value = 1.026702514380507
$ws.Cells.Item(2, 10).Value = 1.047519701752911
$ws.Cells.Item(2, 11).Value = 1.052222093595721
$ws.Cells.Item(2, 12).Value = 1.05231398614712
$ws.Cells.Item(2, 13).Value = 1.060645340779253
$ws.Cells.Item(2, 14).Value = 1.049007299881589

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.04390702417568
$ws.Cells.Item(3, 4).Value = 1.050829055176828
$ws.Cells.Item(3, 5).Value = 1.05090318040054
$ws.Cells.Item(3, 6).Value = 1.059468206855279
$ws.Cells.Item(3, 9).Value = 1.026887458352109
$ws.Cells.Item(3, 10).Value = 1.048626704411167
$ws.Cells.Item(3, 11).Value = 1.053396865821758
$ws.Cells.Item(3, 12).Value = 1.053470799833194
$ws.Cells.Item(3, 13).Value = 1.062013926304268
$ws.Cells.Item(3, 14).Value = 1.050115874610596

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.04484661572813
$ws.Cells.Item(4, 4).Value = 1.051703532563204
$ws.Cells.Item(4, 5).Value = 1.051766047963862
$ws.Cells.Item(4, 6).Value = 1.060462160001686
$ws.Cells.Item(4, 9).Value = 1.026997440722205
$ws.Cells.Item(4, 10).Value = 1.049335499718999
$ws.Cells.Item(4, 11).Value = 1.054148422710619
$ws.Cells.Item(4, 12).Value = 1.054210785064856
$ws.Cells.Item(4, 13).Value = 1.062885796769004
$ws.Cells.Item(4, 14).Value = 1.050825676489064

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.045239930615996
$ws.Cells.Item(5, 4).Value = 1.052069246917554
$ws.Cells.Item(5, 5).Value = 1.052126893946919
$ws.Cells.Item(5, 6).Value = 1.060876877452781
$ws.Cells.Item(5, 9).Value = 1.027041363632053
$ws.Cells.Item(5, 10).Value = 1.049631696560357
$ws.Cells.Item(5, 11).Value = 1.054462337231275
$ws.Cells.Item(5, 12).Value = 1.054519846311849
$ws.Cells.Item(5, 13).Value = 1.06324907775588
$ws.Cells.Item(5, 14).Value = 1.051122293963911

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.045305871415692
$ws.Cells.Item(6, 4).Value = 1.052130540173136
$ws.Cells.Item(6, 5).Value = 1.052187370501709
$ws.Cells.Item(6, 6).Value = 1.060946326944307
$ws.Cells.Item(6, 9).Value = 1.02704860295966
$ws.Cells.Item(6, 10).Value = 1.049681325370345
$ws.Cells.Item(6, 11).Value = 1.054514925757602
$ws.Cells.Item(6, 12).Value = 1.054571620617615
$ws.Cells.Item(6, 13).Value = 1.06330988416365
$ws.Cells.Item(6, 14).Value = 1.051171993252501

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.044851877832406
$ws.Cells.Item(7, 4).Value = 1.051708426758319
$ws.Cells.Item(7, 5).Value = 1.051770877059166
$ws.Cells.Item(7, 6).Value = 1.060467713780345
$ws.Cells.Item(7, 9).Value = 1.026998036705878
$ws.Cells.Item(7, 10).Value = 1.049339464490071
$ws.Cells.Item(7, 11).Value = 1.054152625242893
$ws.Cells.Item(7, 12).Value = 1.054214922702739
$ws.Cells.Item(7, 13).Value = 1.062890663697262
$ws.Cells.Item(7, 14).Value = 1.050829646890566

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.042939867636183
$ws.Cells.Item(8, 4).Value = 1.049927748527896
$ws.Cells.Item(8, 5).Value = 1.050013795259402
$ws.Cells.Item(8, 6).Value = 1.058440463831191
$ws.Cells.Item(8, 9).Value = 1.026767026966958
$ws.Cells.Item(8, 10).Value = 1.047895383946798
$ws.Cells.Item(8, 11).Value = 1.052620905514638
$ws.Cells.Item(8, 12).Value = 1.05270671867564
$ws.Cells.Item(8, 13).Value = 1.061110712345018
$ws.Cells.Item(8, 14).Value = 1.049383515587283

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.039515016698616
$ws.Cells.Item(9, 4).Value = 1.046727086091111
$ws.Cells.Item(9, 5).Value = 1.046855128250006
$ws.Cells.Item(9, 6).Value = 1.054765534017671
$ws.Cells.Item(9, 9).Value = 1.026285444718819
$ws.Cells.Item(9, 10).Value = 1.045292399925421
$ws.Cells.Item(9, 11).Value = 1.049855063281246
$ws.Cells.Item(9, 12).Value = 1.049982697122426
$ws.Cells.Item(9, 13).Value = 1.057868087311537
$ws.Cells.Item(9, 14).Value = 1.046776835030029

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.037192959588801
$ws.Cells.Item(10, 4).Value = 1.044549621222681
$ws.Cells.Item(10, 5).Value = 1.04470595580048
$ws.Cells.Item(10, 6).Value = 1.05224459536142
$ws.Cells.Item(10, 9).Value = 1.025913787246073
$ws.Cells.Item(10, 10).Value = 1.043516644489878
$ws.Cells.Item(10, 11).Value = 1.047964948565619
$ws.Cells.Item(10, 12).Value = 1.048120734699971
$ws.Cells.Item(10, 13).Value = 1.055633140552771
$ws.Cells.Item(10, 14).Value = 1.044998557818086

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.036177939141536
$ws.Cells.Item(11, 4).Value = 1.043596053522936
$ws.Cells.Item(11, 5).Value = 1.043764712841733
$ws.Cells.Item(11, 6).Value = 1.051135690434684
$ws.Cells.Item(11, 9).Value = 1.025740720676983
$ws.Cells.Item(11, 10).Value = 1.042737841059331
$ws.Cells.Item(11, 11).Value = 1.047135223795441
$ws.Cells.Item(11, 12).Value = 1.047303267995196
$ws.Cells.Item(11, 13).Value = 1.054647575697463
$ws.Cells.Item(11, 14).Value = 1.044218648397338

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.035799449524594
$ws.Cells.Item(12, 4).Value = 1.043240215994828
$ws.Cells.Item(12, 5).Value = 1.043413464710395
$ws.Cells.Item(12, 6).Value = 1.050721148065864
$ws.Cells.Item(12, 9).Value = 1.025674599837404
$ws.Cells.Item(12, 10).Value = 1.042447046433167
$ws.Cells.Item(12, 11).Value = 1.046825301161534
$ws.Cells.Item(12, 12).Value = 1.046997909056421
$ws.Cells.Item(12, 13).Value = 1.054278775085686
$ws.Cells.Item(12, 14).Value = 1.04392744080945

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.035880703667571
$ws.Cells.Item(13, 4).Value = 1.043316619078783
$ws.Cells.Item(13, 5).Value = 1.043488882830078
$ws.Cells.Item(13, 6).Value = 1.050810189254916
$ws.Cells.Item(13, 9).Value = 1.025688866292393
$ws.Cells.Item(13, 10).Value = 1.042509491745879
$ws.Cells.Item(13, 11).Value = 1.046891859206927
$ws.Cells.Item(13, 12).Value = 1.047063487693505
$ws.Cells.Item(13, 13).Value = 1.054358007743544
$ws.Cells.Item(13, 14).Value = 1.043989974801668

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.036146683155049
$ws.Cells.Item(14, 4).Value = 1.043566673521652
$ws.Cells.Item(14, 5).Value = 1.043735711966107
$ws.Cells.Item(14, 6).Value = 1.051101478463324
$ws.Cells.Item(14, 9).Value = 1.025735292644082
$ws.Cells.Item(14, 10).Value = 1.042713834894804
$ws.Cells.Item(14, 11).Value = 1.047109640876898
$ws.Cells.Item(14, 12).Value = 1.047278062097084
$ws.Cells.Item(14, 13).Value = 1.054617146213994
$ws.Cells.Item(14, 14).Value = 1.044194608141303

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.036310366726508
$ws.Cells.Item(15, 4).Value = 1.043720522049874
$ws.Cells.Item(15, 5).Value = 1.043887574802229
$ws.Cells.Item(15, 6).Value = 1.051280599443684
$ws.Cells.Item(15, 9).Value = 1.025763653733383
$ws.Cells.Item(15, 10).Value = 1.042839536234584
$ws.Cells.Item(15, 11).Value = 1.047243593749296
$ws.Cells.Item(15, 12).Value = 1.047410040267731
$ws.Cells.Item(15, 13).Value = 1.054776448561862
$ws.Cells.Item(15, 14).Value = 1.044320487991404

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.037260119290489
$ws.Cells.Item(16, 4).Value = 1.044612678105178
$ws.Cells.Item(16, 5).Value = 1.04476819631628
$ws.Cells.Item(16, 6).Value = 1.052317821133348
$ws.Cells.Item(16, 9).Value = 1.025925016331327
$ws.Cells.Item(16, 10).Value = 1.043568120509196
$ws.Cells.Item(16, 11).Value = 1.048019774274153
$ws.Cells.Item(16, 12).Value = 1.04817474833755
$ws.Cells.Item(16, 13).Value = 1.055698170391344
$ws.Cells.Item(16, 14).Value = 1.045050106939255

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.03785329497949
$ws.Cells.Item(17, 4).Value = 1.045169415096865
$ws.Cells.Item(17, 5).Value = 1.04531771802028
$ws.Cells.Item(17, 6).Value = 1.052963775482509
$ws.Cells.Item(17, 9).Value = 1.026022976824906
$ws.Cells.Item(17, 10).Value = 1.044022476111178
$ws.Cells.Item(17, 11).Value = 1.0485036082175
$ws.Cells.Item(17, 12).Value = 1.048651404308034
$ws.Cells.Item(17, 13).Value = 1.056271543482714
$ws.Cells.Item(17, 14).Value = 1.045505107778307

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038198364591761
$ws.Cells.Item(18, 4).Value = 1.045493119069565
$ws.Cells.Item(18, 5).Value = 1.045637220548645
$ws.Cells.Item(18, 6).Value = 1.05333888146671
$ws.Cells.Item(18, 9).Value = 1.026078945463664
$ws.Cells.Item(18, 10).Value = 1.04428654170927
$ws.Cells.Item(18, 11).Value = 1.048784732858904
$ws.Cells.Item(18, 12).Value = 1.048928348693581
$ws.Cells.Item(18, 13).Value = 1.056604266001339
$ws.Cells.Item(18, 14).Value = 1.045769548379837

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.038315869307558
$ws.Cells.Item(19, 4).Value = 1.045603319742117
$ws.Cells.Item(19, 5).Value = 1.045745989833537
$ws.Cells.Item(19, 6).Value = 1.053466501222861
$ws.Cells.Item(19, 9).Value = 1.026097831232171
$ws.Cells.Item(19, 10).Value = 1.044376420504082
$ws.Cells.Item(19, 11).Value = 1.048880405499812
$ws.Cells.Item(19, 12).Value = 1.049022597074765
$ws.Cells.Item(19, 13).Value = 1.056717425871827
$ws.Cells.Item(19, 14).Value = 1.045859554812848

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.037789748189684
$ws.Cells.Item(20, 4).Value = 1.045109789376692
$ws.Cells.Item(20, 5).Value = 1.045258865694314
$ws.Cells.Item(20, 6).Value = 1.052894643556398
$ws.Cells.Item(20, 9).Value = 1.026012587705531
$ws.Cells.Item(20, 10).Value = 1.043973826747977
$ws.Cells.Item(20, 11).Value = 1.048451810108494
$ws.Cells.Item(20, 12).Value = 1.048600375649595
$ws.Cells.Item(20, 13).Value = 1.056210203791567
$ws.Cells.Item(20, 14).Value = 1.04545638932743

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.036068399540243
$ws.Cells.Item(21, 4).Value = 1.043493084248277
$ws.Cells.Item(21, 5).Value = 1.043663072146718
$ws.Cells.Item(21, 6).Value = 1.051015774444674
$ws.Cells.Item(21, 9).Value = 1.025721672037495
$ws.Cells.Item(21, 10).Value = 1.042653702898652
$ws.Cells.Item(21, 11).Value = 1.047045557489447
$ws.Cells.Item(21, 12).Value = 1.047214922877252
$ws.Cells.Item(21, 13).Value = 1.054540911761518
$ws.Cells.Item(21, 14).Value = 1.044134390750819

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.034977622440786
$ws.Cells.Item(22, 4).Value = 1.042467092531988
$ws.Cells.Item(22, 5).Value = 1.04265029472673
$ws.Cells.Item(22, 6).Value = 1.049819124577815
$ws.Cells.Item(22, 9).Value = 1.025528130931748
$ws.Cells.Item(22, 10).Value = 1.041814924520298
$ws.Cells.Item(22, 11).Value = 1.046151389481651
$ws.Cells.Item(22, 12).Value = 1.046333893485316
$ws.Cells.Item(22, 13).Value = 1.053475614568149
$ws.Cells.Item(22, 14).Value = 1.043294421210959

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.035556679932524
$ws.Cells.Item(23, 4).Value = 1.043011901800644
$ws.Cells.Item(23, 5).Value = 1.043188092443821
$ws.Cells.Item(23, 6).Value = 1.050454959483391
$ws.Cells.Item(23, 9).Value = 1.025631742986492
$ws.Cells.Item(23, 10).Value = 1.042260416752527
$ws.Cells.Item(23, 11).Value = 1.046626363030159
$ws.Cells.Item(23, 12).Value = 1.046801896127908
$ws.Cells.Item(23, 13).Value = 1.054041855650669
$ws.Cells.Item(23, 14).Value = 1.043740546093257

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.037818465089422
$ws.Cells.Item(24, 4).Value = 1.045136734858522
$ws.Cells.Item(24, 5).Value = 1.045285461690633
$ws.Cells.Item(24, 6).Value = 1.052925886451072
$ws.Cells.Item(24, 9).Value = 1.026017285716356
$ws.Cells.Item(24, 10).Value = 1.043995812242738
$ws.Cells.Item(24, 11).Value = 1.048475218803507
$ws.Cells.Item(24, 12).Value = 1.048623436642837
$ws.Cells.Item(24, 13).Value = 1.056237925859222
$ws.Cells.Item(24, 14).Value = 1.045478406044116

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.040407157614159
$ws.Cells.Item(25, 4).Value = 1.047562126109793
$ws.Cells.Item(25, 5).Value = 1.047679260581157
$ws.Cells.Item(25, 6).Value = 1.055727951611693
$ws.Cells.Item(25, 9).Value = 1.026418823642792
$ws.Cells.Item(25, 10).Value = 1.045972366797057
$ws.Cells.Item(25, 11).Value = 1.05057814345933
$ws.Cells.Item(25, 12).Value = 1.05069491882206
$ws.Cells.Item(25, 13).Value = 1.058719142409143
$ws.Cells.Item(25, 14).Value = 1.047457767532616
